$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 05:21"

# Row 8 - Peru
$ws.Range("B8").Value = 716670
$ws.Range("D8").Value = 552214
$ws.Range("E8").Value = 133986
$ws.Range("H8").Value = 30470

# Row 20 - Pakistan
$ws.Range("B20").Value = 300955
$ws.Range("C20").Value = 584
$ws.Range("D20").Value = 288536
$ws.Range("E20").Value = 6046
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 6373

# Row 33 - Kazajistan
$ws.Range("B33").Value = 106729
$ws.Range("C33").Value = 68
$ws.Range("E33").Value = 4686

# Row 39 - Belgica
$ws.Range("B39").Value = 91537
$ws.Range("C39").Value = 969
$ws.Range("D39").Value = 18689
$ws.Range("E39").Value = 62929
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 9919

# Row 160 - Belice
$ws.Range("B160").Value = 1435
$ws.Range("C160").Value = 36
$ws.Range("D160").Value = 404
$ws.Range("E160").Value = 1012

# Row 172 - Islas Turcas y Caicos
$ws.Range("B172").Value = 641
$ws.Range("C172").Value = 3
$ws.Range("E172").Value = 366

# Row 173 - San Martin (Parte Holandesa)
$ws.Range("B173").Value = 531
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 387
$ws.Range("E173").Value = 125

# Row 198 - San Vicente y las Granadinas
$ws.Range("B198").Value = 64
$ws.Range("C198").Value = 2
$ws.Range("E198").Value = 3
